$wb = $excel.ActiveWorkbook

# --- Update the "Conversion del dia" text on sheet "Hoja1" ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 13.52 = 54825.63 pesos`n✅ 54825.63 pesos = 13.45 = 974.95 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update the rate cells on sheet "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 73.98
$wsTasas.Range("O10").Value = 4056
$wsTasas.Range("N12").Value = 4077
$wsTasas.Range("O12").Value = 72.5
